$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the "total" summary row from row 4 up to row 3, closing the empty
# gap left by the (now removed) rows 2-3 — fixes the extra-time sum row.
$cols = @("A", "G", "H", "I")
foreach ($col in $cols) {
    $src = $ws.Range($col + "4")
    $dst = $ws.Range($col + "3")
    $dst.Value = $src.Text
    $dst.Font.Color = $src.Font.Color
    $dst.NumberFormat = $src.NumberFormat
}

$ws.Range("A4:I4").Delete()
